$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '27.304.32'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '1.863.88'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +1.24%  '
Set-TextValue $ws.Range('D5') '313.31'
$ws.Range('E6').Value = '  +1.22%  '
Set-TextValue $ws.Range('D7') '0.4807'
$ws.Range('E7').Value = '  +2.31%  '
Set-TextValue $ws.Range('D8') '0.3730'
$ws.Range('E8').Value = '  +2.15%  '
Set-TextValue $ws.Range('D9') '0.07462'
$ws.Range('E9').Value = '  +4.53%  '
Set-TextValue $ws.Range('D10') '0.9379'
$ws.Range('E10').Value = '  +2.16%  '
Set-TextValue $ws.Range('D11') '20.70'
$ws.Range('E11').Value = '  +6.05%  '
Set-TextValue $ws.Range('D12') '0.07880'
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('D13').Value = '1.861.19'
$ws.Range('E13').Value = '  +0.41%  '
$ws.Range('E14').Value = '  +2.91%  '
Set-TextValue $ws.Range('D15') '6.540'
$ws.Range('E15').Value = '  +2.25%  '
Set-TextValue $ws.Range('D16') '90.41'
$ws.Range('E16').Value = '  +2.73%  '
Set-TextValue $ws.Range('D17') '1.021'
$ws.Range('E17').Value = '  +1.16%  '
Set-TextValue $ws.Range('D18') '0.000008811'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('E19').Value = '  +1.21%  '
Set-TextValue $ws.Range('D20') '14.84'
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('D21').Value = '27.334.54'
$ws.Range('E21').Value = '  +1.61%  '
Set-TextValue $ws.Range('D22') '5.127'
Set-TextValue $ws.Range('D23') '10.70'
$ws.Range('E23').Value = '  +0.86%  '
Set-TextValue $ws.Range('D24') '1.960'
$ws.Range('E24').Value = '  +2.11%  '
Set-TextValue $ws.Range('D25') '154.24'
$ws.Range('E26').Value = '  +2.13%  '
Set-TextValue $ws.Range('D27') '2.011'
$ws.Range('E27').Value = '  +0.08%  '
Set-TextValue $ws.Range('D28') '116.05'
$ws.Range('E28').Value = '  +1.63%  '
Set-TextValue $ws.Range('D29') '5.003'
$ws.Range('E29').Value = '  +2.88%  '
Set-TextValue $ws.Range('D30') '0.08919'
$ws.Range('E30').Value = '  +1.18%  '
Set-TextValue $ws.Range('D31') '3.349'
$ws.Range('E31').Value = '  +4.04%  '
$ws.Range('E32').Value = '  +2.17%  '
Set-TextValue $ws.Range('D33') '4.577'
$ws.Range('E33').Value = '  +2.35%  '
Set-TextValue $ws.Range('D34') '0.7439'
$ws.Range('E34').Value = '  -0.06%  '
Set-TextValue $ws.Range('D35') '2.674'
$ws.Range('E35').Value = '  -2.92%  '
Set-TextValue $ws.Range('D36') '0.02055'
$ws.Range('E36').Value = '  +5.88%  '
Set-TextValue $ws.Range('D37') '1.125'
$ws.Range('E37').Value = '  +3.43%  '
Set-TextValue $ws.Range('D38') '0.05294'
$ws.Range('E38').Value = '  +1.60%  '
Set-TextValue $ws.Range('D39') '0.5376'
$ws.Range('E39').Value = '  +3.62%  '
Set-TextValue $ws.Range('D40') '7.137'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('E41').Value = '  +1.77%  '
Set-TextValue $ws.Range('D42') '8.393'
$ws.Range('E42').Value = '  +2.97%  '
Set-TextValue $ws.Range('D43') '10.64'
$ws.Range('E43').Value = '  +1.55%  '
Set-TextValue $ws.Range('D44') '0.4836'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('E45').Value = '  +1.29%  '
Set-TextValue $ws.Range('D46') '1.676'
$ws.Range('E46').Value = '  +5.29%  '
Set-TextValue $ws.Range('D47') '103.24'
$ws.Range('E47').Value = '  +1.50%  '
Set-TextValue $ws.Range('D48') '66.73'
$ws.Range('E48').Value = '  +2.42%  '
Set-TextValue $ws.Range('D49') '0.06090'
$ws.Range('E49').Value = '  +0.98%  '
Set-TextValue $ws.Range('D50') '0.9011'
$ws.Range('E50').Value = '  +1.90%  '
Set-TextValue $ws.Range('D51') '36.86'
$ws.Range('E51').Value = '  +1.76%  '
